$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 697:698, pushing all following rows down by two
# (old row 697 "2026/12/29" becomes new row 699, ..., old row 738 becomes new row 740)
$ws.Range("A697:A698").EntireRow.Insert()

# New row 697: 2026/01/25, 日, 23, 17
# Leading apostrophe forces the date-looking text to stay literal text
# (matches the existing inlineStr cells used throughout column A/B).
$ws.Range("A697").Value = "'2026/01/25"
$ws.Range("B697").Value = "日"
$ws.Range("C697").Value = 23
$ws.Range("D697").Value = 17

# New row 698: 2026/01/26, 月, 2, 18
$ws.Range("A698").Value = "'2026/01/26"
$ws.Range("B698").Value = "月"
$ws.Range("C698").Value = 2
$ws.Range("D698").Value = 18
